# Add a "Save" column (H) to the s_vals sheet, matching header style of
# the existing header row and filling in numeric save values per row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style from the existing "sum" header (G1) onto the new H1 header
# cell, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill in the new "Save" values for each data row.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 0
